# Updates the "Estado de Cuenta" worker/period table (rows 16-37) so the
# data is grouped by period (1607 first, then 1608) with two extra rows for
# periods 1610 and 1611, matching the refreshed EC database export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# r=row => (DocNumber, Name, Period)
$rows = @{
    16 = @("45693199",   "JANIA DEL CARMEN CASTRO SIBAJA",       "1607")
    17 = @("3811649",    "AMIR JOSE HERNANDEZ HERNANDEZ",        "1607")
    18 = @("32291589",   "LUZ AMPARO RIVAS LOPEZ",                "1607")
    19 = @("45488572",   "RUTH DEL CARMEN GASTELBONDO HERRERA",  "1607")
    20 = @("18002992",   "OSVALDO ANTONIO YEPES DE ORO",          "1607")
    21 = @("45592205",   "ESTILITA SIMANCAS REBOLLEDO",           "1607")
    22 = @("45545576",   "ANA ISABEL OVIEDO CABRERA",             "1607")
    23 = @("1140866934", "JESSICA NARDELLY GUARNIZO RAMOS",       "1607")
    24 = @("1047428222", "VANESA PAOLA FERNANDEZ HERNANDEZ",      "1607")
    25 = @("1047373116", "SINDY YANINA LOPEZ RIVERA",             "1607")
    26 = @("45693199",   "JANIA DEL CARMEN CASTRO SIBAJA",        "1608")
    27 = @("3811649",    "AMIR JOSE HERNANDEZ HERNANDEZ",         "1608")
    28 = @("32291589",   "LUZ AMPARO RIVAS LOPEZ",                "1608")
    29 = @("45488572",   "RUTH DEL CARMEN GASTELBONDO HERRERA",   "1608")
    30 = @("18002992",   "OSVALDO ANTONIO YEPES DE ORO",          "1608")
    31 = @("45592205",   "ESTILITA SIMANCAS REBOLLEDO",           "1608")
    32 = @("45545576",   "ANA ISABEL OVIEDO CABRERA",             "1608")
    33 = @("1140866934", "JESSICA NARDELLY GUARNIZO RAMOS",       "1608")
    34 = @("1047428222", "VANESA PAOLA FERNANDEZ HERNANDEZ",      "1608")
    35 = @("1047373116", "SINDY YANINA LOPEZ RIVERA",             "1608")
    36 = @("18002992",   "OSVALDO ANTONIO YEPES DE ORO",          "1610")
    37 = @("18002992",   "OSVALDO ANTONIO YEPES DE ORO",          "1611")
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 3).Value = $vals[0]  # C: N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $vals[1]  # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $vals[2]  # E: Periodo Mora
}
